# fix: response member dan sub kategori untuk invoice
#
# Adds a new "type payment" column (D) to the demoSubCategory sheet, classifying
# each existing sub-category row as a "once" / "monthly" / "anytime" payment type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell: D1 = "type payment" -----------------------------------
# Match the look of the existing header row (A1:C1): bold font on a yellow fill.
$ws.Range("D1").Value = "type payment"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").Interior.ColorIndex = $ws.Range("A1").Interior.ColorIndex

# --- New data column: D2:D12 = payment cadence for each row -----------------
$paymentTypes = @(
    "once",     # simpanan pokok
    "monthly",  # simpanan wajib
    "monthly",  # simpanan wajib khusus
    "anytime",  # simpanan sukarela
    "anytime",  # tabungan rekreasi
    "monthly",  # piutang s/p
    "monthly",  # piutang dagang
    "anytime",  # pembelian barang
    "anytime",  # penjualan barang
    "monthly",  # pembayaran angsuran
    "monthly"   # jasa s/p
)

for ($i = 0; $i -lt $paymentTypes.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = $paymentTypes[$i]
    $cell.Font.Bold = $false
}

# --- Column widths ------------------------------------------------------------
# Column D now holds the longer "type payment" values, so widen it; column E
# (previously sharing a width with D) keeps the original narrow width.
$ws.Range("D1").EntireColumn.ColumnWidth = 16.75
$ws.Range("E1").EntireColumn.ColumnWidth = 7.83
